# Weekly fruit/vegetable price update:
# Insert one new daily record for "Zapallo italiano" at row 46 (pushing all
# subsequent rows down by one), matching the new sheet dimension A1:R142.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before the current row 46, shifting rows 46:141 down to 47:142.
$ws.Rows.Item(46).Insert()

# Populate the newly inserted row 46 with the latest price record.
$ws.Cells.Item(46, 1).Value = 11
$ws.Cells.Item(46, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(46, 3).Value = "Bíobío"
$ws.Cells.Item(46, 4).Value = 44775
$ws.Cells.Item(46, 5).Value = 8
$ws.Cells.Item(46, 6).Value = 100112032
$ws.Cells.Item(46, 7).Value = "Zapallo italiano"
$ws.Cells.Item(46, 8).Value = "Sin especificar"
$ws.Cells.Item(46, 9).Value = "Primera"
$ws.Cells.Item(46, 10).Value = 180
$ws.Cells.Item(46, 11).Value = 17000
$ws.Cells.Item(46, 12).Value = 18000
$ws.Cells.Item(46, 13).Value = 17556
$ws.Cells.Item(46, 14).Value = "$/caja 50 unidades"
$ws.Cells.Item(46, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(46, 16).Value = 351
$ws.Cells.Item(46, 17).Value = 50
$ws.Cells.Item(46, 18).Value = "Hortaliza"
